# Update the StructureDefinition-claim-adjustment-type workbook:
# - rebrand from ibm.com/Alvearie Team to linuxforhealth.org/LinuxForHealth Team
# - bump version 7.0.0 -> 8.0.0
# - update publish date
# - update the payer-claim-adjustment-type ValueSet URL (also drops "wh-" prefix)

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-adjustment-type"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-adjustment-type"
$elements.Range("Y7").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/payer-claim-adjustment-type"
